$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.931.27'
$ws.Range('E2').Value = '  -0.59%  '

$ws.Range('D3').Value = '1.952.03'
$ws.Range('E3').Value = '  -0.89%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').Value = "'" + '242.76'
$ws.Range('E5').Value = '  -2.42%  '

$ws.Range('D6').Value = "'" + '1.002'
$ws.Range('E6').Value = '  -0.11%  '

$ws.Range('D7').Value = "'" + '0.4865'
$ws.Range('E7').Value = '  -0.34%  '

$ws.Range('D8').Value = "'" + '0.2932'
$ws.Range('E8').Value = '  -1.04%  '

$ws.Range('D9').Value = "'" + '0.07007'
$ws.Range('E9').Value = '  +2.55%  '

$ws.Range('D10').Value = "'" + '19.54'
$ws.Range('E10').Value = '  +1.47%  '

$ws.Range('D11').Value = "'" + '107.02'
$ws.Range('E11').Value = '  -0.79%  '

$ws.Range('D12').Value = '1.955.34'
$ws.Range('E12').Value = '  -0.88%  '

$ws.Range('D13').Value = "'" + '0.07754'
$ws.Range('E13').Value = '  -0.55%  '

$ws.Range('D14').Value = "'" + '5.344'
$ws.Range('E14').Value = '  -2.11%  '

$ws.Range('D15').Value = "'" + '0.6984'
$ws.Range('E15').Value = '  -0.94%  '

$ws.Range('D16').Value = "'" + '277.55'
$ws.Range('E16').Value = '  -3.64%  '

$ws.Range('D17').Value = '30.950.78'
$ws.Range('E17').Value = '  -0.56%  '

$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').Value = "'" + '13.19'
$ws.Range('E18').Value = '  -0.47%  '

$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = "'" + '0.000007727'
$ws.Range('E19').Value = '  -0.44%  '

$ws.Range('D20').Value = '2.209.80'
$ws.Range('E20').Value = '  -0.64%  '

$ws.Range('E21').Value = '  -0.08%  '

$ws.Range('D22').Value = "'" + '5.479'
$ws.Range('E22').Value = '  -2.82%  '

$ws.Range('D23').Value = "'" + '1.001'
$ws.Range('E23').Value = '  -0.70%  '

$ws.Range('D24').Value = "'" + '6.475'
$ws.Range('E24').Value = '  -2.44%  '

$ws.Range('D25').Value = "'" + '9.728'
$ws.Range('E25').Value = '  -3.09%  '

$ws.Range('D26').Value = "'" + '168.86'
$ws.Range('E26').Value = '  -1.08%  '

$ws.Range('D27').Value = "'" + '19.63'
$ws.Range('E27').Value = '  -2.20%  '

$ws.Range('D28').Value = "'" + '2.163'
$ws.Range('E28').Value = '  -1.54%  '

$ws.Range('D29').Value = "'" + '0.1043'
$ws.Range('E29').Value = '  -2.54%  '

$ws.Range('D30').Value = "'" + '1.396'
$ws.Range('E30').Value = '  -3.45%  '

$ws.Range('D31').Value = "'" + '4.632'
$ws.Range('E31').Value = '  -4.39%  '

$ws.Range('D32').Value = "'" + '1.559'
$ws.Range('E32').Value = '  -2.95%  '

$ws.Range('D33').Value = "'" + '4.393'
$ws.Range('E33').Value = '  -3.00%  '

$ws.Range('D34').Value = "'" + '0.04874'
$ws.Range('E34').Value = '  -4.63%  '

$ws.Range('D35').Value = "'" + '0.7513'
$ws.Range('E35').Value = '  -2.86%  '

$ws.Range('E36').Value = '  -0.88%  '

$ws.Range('D37').Value = "'" + '2.733'
$ws.Range('E37').Value = '  -0.15%  '

$ws.Range('D38').Value = "'" + '0.01990'
$ws.Range('E38').Value = '  -3.04%  '

$ws.Range('D39').Value = "'" + '2.677'
$ws.Range('E39').Value = '  -2.01%  '

$ws.Range('D40').Value = "'" + '6.518'
$ws.Range('E40').Value = '  +0.06%  '

$ws.Range('D41').Value = "'" + '77.88'
$ws.Range('E41').Value = '  +6.66%  '

$ws.Range('D42').Value = "'" + '2.095'
$ws.Range('E42').Value = '  -2.07%  '

$ws.Range('D43').Value = "'" + '0.8940'
$ws.Range('E43').Value = '  +0.32%  '

$ws.Range('D44').Value = "'" + '109.08'
$ws.Range('E44').Value = '  -0.98%  '

$ws.Range('D45').Value = "'" + '0.4430'
$ws.Range('E45').Value = '  -1.39%  '

$ws.Range('D46').Value = "'" + '0.9999'
$ws.Range('E46').Value = '  -0.31%  '

$ws.Range('D47').Value = "'" + '7.757'
$ws.Range('E47').Value = '  +2.81%  '

$ws.Range('D48').Value = "'" + '989.27'
$ws.Range('E48').Value = '  -0.84%  '

$ws.Range('D49').Value = "'" + '0.1247'
$ws.Range('E49').Value = '  -1.82%  '

$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = "'" + '35.90'
$ws.Range('E50').Value = '  -0.62%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'" + '9.187'
$ws.Range('E51').Value = '  -3.30%  '

